$d = $word.ActiveDocument
Write-Output $d.Paragraphs.Count
